$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new commit-log row (row 14) documenting the latest commit:
#   2025-08-11 20:00 | e7098fd | Fix category card weekly breakdown arrow
#   visibility issue | 2 files changed | +17 / -5 | session description | Local
$ws.Range("A14").Value = "2025-08-11 20:00"
$ws.Range("B14").Value = "e7098fd"
$ws.Range("C14").Value = "Fix category card weekly breakdown arrow visibility issue"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 17
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = "Fixed category card weekly breakdown arrow visibility CSS logic"
$ws.Range("H14").Value = "Local"
